$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.664.94"
$ws.Range("E2").Value = "  -0.40%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.263.32"
$ws.Range("E3").Value = "  -0.50%  "

# Row 4
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.78"
$ws.Range("E5").Value = "  +0.06%  "

# Row 6
$ws.Range("E6").Value = "  +0.56%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "76.93"
$ws.Range("E7").Value = "  +6.55%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.640"
$ws.Range("E9").Value = "  -3.11%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.71"
$ws.Range("E10").Value = "  +3.35%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0970"
$ws.Range("E11").Value = "  -0.07%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.25"
$ws.Range("E12").Value = "  -2.36%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.106"
$ws.Range("E13").Value = "  +1.80%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.600.71"
$ws.Range("E14").Value = "  -0.55%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.01"
$ws.Range("E15").Value = "  +1.43%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.863"
$ws.Range("E16").Value = "  -2.66%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.267.35"
$ws.Range("E17").Value = "  -0.57%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.530.27"
$ws.Range("E18").Value = "  -0.57%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0987"
$ws.Range("E19").Value = "  -1.54%  "

# Row 20
$ws.Range("E20").Value = "  -2.04%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.02"
$ws.Range("E21").Value = "  -1.53%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "235.62"
$ws.Range("E22").Value = "  +0.02%  "

# Row 23
$ws.Range("E23").Value = "  -0.01%  "

# Row 24
$ws.Range("E24").Value = "  +0.07%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.75"
$ws.Range("E25").Value = "  -7.27%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.31"
$ws.Range("E26").Value = "  -0.50%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.38"
$ws.Range("E27").Value = "  -2.32%  "

# Row 28
$ws.Range("E28").Value = "  -1.03%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.40"
$ws.Range("E29").Value = "  -0.06%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.91"
$ws.Range("E30").Value = "  -0.53%  "

# Row 31
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0860"
$ws.Range("E31").Value = "  +6.48%  "

# Row 32
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.39"
$ws.Range("E32").Value = "  -1.51%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.123"
$ws.Range("E33").Value = "  -3.21%  "

# Row 34
$ws.Range("E34").Value = "  -1.68%  "

# Row 35
$ws.Range("E35").Value = "  +0.96%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.56"
$ws.Range("E36").Value = "  +1.73%  "

# Row 37
$ws.Range("E37").Value = "  -0.85%  "

# Row 38
$ws.Range("E38").Value = "  -3.13%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.90"
$ws.Range("E39").Value = "  +8.72%  "

# Row 40
$ws.Range("E40").Value = "  -3.05%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.86"
$ws.Range("E41").Value = "  +0.70%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.207"
$ws.Range("E42").Value = "  +0.02%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "61.43"
$ws.Range("E43").Value = "  -1.01%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "108.10"
$ws.Range("E44").Value = "  +13.48%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.85"
$ws.Range("E45").Value = "  -4.90%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.66"
$ws.Range("E46").Value = "  -7.62%  "

# Row 47
$ws.Range("E47").Value = "  -1.84%  "

# Row 48
$ws.Range("E48").Value = "  -0.33%  "

# Row 49
$ws.Range("E49").Value = "  -2.09%  "

# Row 50
$ws.Range("E50").Value = "  -2.63%  "

# Row 51
$ws.Range("B51").Value = "Bonk"
$ws.Range("C51").Value = "https://coinranking.com/coin/jCd_nuYCH+bonk-bonk"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0000333"
$ws.Range("E51").Value = "  +129.21%  "
